# Fix minor grammar issues in the report.
#
# 1. "...only yields a revenue of $6,840." -> "...only yields revenue of $6,840."
# 2. "...Expected revenues (calculated..." -> "...Expected revenue (calculated..."  (Targeted Addresses section)
# 3. "...the profits -$97,635..." -> "...the profit is -$97,635..." (Recommendations section)
# 4. "...expected revenues with a target population..." -> "...expected revenue with a target population..." (Recommendations section)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "targeting everyone only yields a revenue of `$6,840.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "targeting everyone only yields revenue of `$6,840.",
    2)

$d.Content.Find.Execute(
    "provided as a csv file. Expected revenues (calculated",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "provided as a csv file. Expected revenue (calculated",
    2)

$d.Content.Find.Execute(
    "30,000 people, the profits -`$97,635. Similarly, expected revenues with a target population",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "30,000 people, the profit is -`$97,635. Similarly, expected revenue with a target population",
    2)
